$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1004.8333
$ws.Range("I12").Value = 805.4
$ws.Range("J12").Value = 2002
$ws.Range("K12").Value = 805.4
$ws.Range("L12").Value = 2002
$ws.Range("M12").Value = -635.4
$ws.Range("N12").Value = -2342

$ws.Range("H33").Value = 356.16666
$ws.Range("I33").Value = 356.16666
$ws.Range("K33").Value = 356.16666
$ws.Range("M33").Value = -127.16666

$ws.Range("H43").Value = 6650
$ws.Range("I43").Value = 6308.3335
$ws.Range("K43").Value = 6308.3335
$ws.Range("M43").Value = -6239.3335

$ws.Range("H64").Value = 7571.2856
$ws.Range("I64").Value = 5749.75
$ws.Range("K64").Value = 5749.75
$ws.Range("M64").Value = -5501.75

$ws.Range("H67").Value = 7571.2856
$ws.Range("I67").Value = 5749.75
$ws.Range("K67").Value = 5749.75
$ws.Range("M67").Value = -4891.75

$ws.Range("H138").Value = 6430.0386
$ws.Range("J138").Value = 6804.3
$ws.Range("L138").Value = 20412.9
$ws.Range("N138").Value = -30692.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2031.6666
$ws.Range("I5").Value = 1047.5
$ws.Range("J5").Value = 4000
$ws.Range("K5").Value = 1047.5
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = -935.5
$ws.Range("N5").Value = -4224

$ws.Range("H14").Value = 334466.66
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H63").Value = 10250.75
$ws.Range("I63").Value = 5003
$ws.Range("K63").Value = 5003
$ws.Range("M63").Value = -4317

$ws.Range("H66").Value = 10250.75
$ws.Range("I66").Value = 5003
$ws.Range("K66").Value = 25015
$ws.Range("M66").Value = -21583

$ws.Range("H132").Value = 2552.4
$ws.Range("I132").Value = 1940.5
$ws.Range("K132").Value = 5821.5
$ws.Range("M132").Value = -3291.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2031.6666
$ws.Range("I4").Value = 1047.5
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 1047.5
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = -932.5
$ws.Range("N4").Value = -4230

$ws.Range("H15").Value = 48333.332
$ws.Range("J15").Value = 48333.332
$ws.Range("L15").Value = 48333.332
$ws.Range("N15").Value = -48787.332

$ws.Range("H19").Value = 46666.668
$ws.Range("J19").Value = 46666.668
$ws.Range("L19").Value = 46666.668
$ws.Range("N19").Value = -47012.668

$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2377
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11884
$ws.Range("N89").ClearContents()

$ws.Range("H105").Value = 2251.25
$ws.Range("I105").Value = 1001.6667
$ws.Range("K105").Value = 1001.6667
$ws.Range("M105").Value = 745.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 111889.664
$ws.Range("I4").Value = 902.5
$ws.Range("J4").Value = 200679.4
$ws.Range("K4").Value = 902.5
$ws.Range("L4").Value = 200679.4
$ws.Range("M4").Value = -790.5
$ws.Range("N4").Value = -200903.4

$ws.Range("H7").Value = 196.36363
$ws.Range("J7").Value = 394
$ws.Range("L7").Value = 394
$ws.Range("N7").Value = -620

$ws.Range("H15").Value = 1107.4073
$ws.Range("I15").Value = 517.3913
$ws.Range("K15").Value = 517.3913
$ws.Range("M15").Value = -347.3913

$ws.Range("H31").Value = 7506
$ws.Range("I31").Value = 5012
$ws.Range("K31").Value = 5012
$ws.Range("M31").Value = -4717

$ws.Range("H34").Value = 7506
$ws.Range("I34").Value = 5012
$ws.Range("K34").Value = 5012
$ws.Range("M34").Value = -4810

$ws.Range("H69").Value = 37692.08
$ws.Range("J69").Value = 81499.25
$ws.Range("L69").Value = 81499.25
$ws.Range("N69").Value = -82997.25

$ws.Range("H72").Value = 37692.08
$ws.Range("J72").Value = 81499.25
$ws.Range("L72").Value = 244497.75
$ws.Range("N72").Value = -251985.75

$ws.Range("H107").Value = 937
$ws.Range("I107").Value = 628.6667
$ws.Range("J107").Value = 1399.5
$ws.Range("K107").Value = 628.6667
$ws.Range("L107").Value = 1399.5
$ws.Range("M107").Value = 1291.3333
$ws.Range("N107").Value = -5239.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1535.8064
$ws.Range("I4").Value = 1305.4546
$ws.Range("J4").Value = 2098.889
$ws.Range("K4").Value = 3916.3638
$ws.Range("L4").Value = 6296.667
$ws.Range("M4").Value = -3804.3638
$ws.Range("N4").Value = -6520.667

$ws.Range("H12").Value = 81.92308
$ws.Range("I12").Value = 8.333333
$ws.Range("K12").Value = 24.999999
$ws.Range("M12").Value = 148.000001

$ws.Range("H74").Value = 750
$ws.Range("I74").Value = 750
$ws.Range("K74").Value = 2250
$ws.Range("M74").Value = -1189

$ws.Range("H77").Value = 750
$ws.Range("I77").Value = 750
$ws.Range("K77").Value = 6750
$ws.Range("M77").Value = -1446

$ws.Range("H137").Value = 1370
$ws.Range("I137").Value = 1370
$ws.Range("K137").Value = 4110
$ws.Range("M137").Value = 990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4003
$ws.Range("J80").Value = 4006
$ws.Range("L80").Value = 4006
$ws.Range("N80").Value = -6002

$ws.Range("H83").Value = 4003
$ws.Range("J83").Value = 4006
$ws.Range("L83").Value = 20030
$ws.Range("N83").Value = -30014

$ws.Range("H132").Value = 3592.2
$ws.Range("I132").Value = 3592.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10776.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8246.599999999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3003
$ws.Range("I7").Value = 3002
$ws.Range("K7").Value = 3002
$ws.Range("M7").Value = -2890

$ws.Range("H22").Value = 1583.3334
$ws.Range("I22").Value = 1583.3334
$ws.Range("K22").Value = 1583.3334
$ws.Range("M22").Value = -1288.3334

$ws.Range("H27").Value = 1583.3334
$ws.Range("I27").Value = 1583.3334
$ws.Range("K27").Value = 1583.3334
$ws.Range("M27").Value = -1476.3334

$ws.Range("H68").Value = 4249.75
$ws.Range("I68").Value = 4249.75
$ws.Range("K68").Value = 4249.75
$ws.Range("M68").Value = -3500.75

$ws.Range("H71").Value = 4249.75
$ws.Range("I71").Value = 4249.75
$ws.Range("K71").Value = 21248.75
$ws.Range("M71").Value = -17504.75

$ws.Range("H74").Value = 70666.664
$ws.Range("I74").Value = 50000
$ws.Range("J74").Value = 81000
$ws.Range("K74").Value = 50000
$ws.Range("L74").Value = 81000
$ws.Range("M74").Value = -49002
$ws.Range("N74").Value = -82996

$ws.Range("H77").Value = 70666.664
$ws.Range("I77").Value = 50000
$ws.Range("J77").Value = 81000
$ws.Range("K77").Value = 150000
$ws.Range("L77").Value = 243000
$ws.Range("M77").Value = -145008
$ws.Range("N77").Value = -252984

$ws.Range("H82").Value = 7450
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 9900
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 9900
$ws.Range("M82").Value = -4639
$ws.Range("N82").Value = -10622

$ws.Range("H85").Value = 7450
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 9900
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 9900
$ws.Range("M85").Value = -3752
$ws.Range("N85").Value = -12396

$ws.Range("H126").Value = 3003
$ws.Range("I126").Value = 3002
$ws.Range("K126").Value = 9006
$ws.Range("M126").Value = -6536

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 21999
$ws.Range("I21").Value = 21999
$ws.Range("K21").Value = 21999
$ws.Range("M21").Value = -21764

$ws.Range("H35").Value = 21999
$ws.Range("I35").Value = 21999
$ws.Range("K35").Value = 21999
$ws.Range("M35").Value = -21709

$ws.Range("H62").Value = 6664.5
$ws.Range("J62").Value = 8992.5
$ws.Range("L62").Value = 8992.5
$ws.Range("N62").Value = -10240.5

$ws.Range("H64").Value = 10394.5
$ws.Range("J64").Value = 10394.5
$ws.Range("L64").Value = 10394.5
$ws.Range("N64").Value = -10890.5

$ws.Range("H65").Value = 6664.5
$ws.Range("J65").Value = 8992.5
$ws.Range("L65").Value = 44962.5
$ws.Range("N65").Value = -51202.5

$ws.Range("H67").Value = 10394.5
$ws.Range("J67").Value = 10394.5
$ws.Range("L67").Value = 10394.5
$ws.Range("N67").Value = -12110.5

$ws.Range("H105").Value = 39181.547
$ws.Range("J105").Value = 39181.547
$ws.Range("L105").Value = 39181.547
$ws.Range("N105").Value = -46169.547

$ws.Range("H132").Value = 128.33333
$ws.Range("I132").Value = 128.33333
$ws.Range("K132").Value = 384.99999
$ws.Range("M132").Value = 2145.00001
